$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -137.11378193317705
$ws.Range("C2").Value = 0.052560452282368476
$ws.Range("D2").Value = 301.991214785

$ws.Range("B3").Value = -135.31516283528447
$ws.Range("C3").Value = 0.09800251867935554
$ws.Range("D3").Value = 384.322671638

$ws.Range("B4").Value = -135.95637869657497
$ws.Range("C4").Value = 0.09811297264043536
$ws.Range("D4").Value = 448.975615579

$ws.Range("B5").Value = -135.5850539924898
$ws.Range("C5").Value = 0.09669560631741587
$ws.Range("D5").Value = 210.674142406

$ws.Range("B6").Value = -135.09244017942407
$ws.Range("C6").Value = 0.07315644398569157
$ws.Range("D6").Value = 267.443320047

$ws.Range("B7").Value = -134.7311225764285
$ws.Range("C7").Value = 0.09963374805368759
$ws.Range("D7").Value = 306.927288725

$ws.Range("B8").Value = -134.05690452315747
$ws.Range("C8").Value = 0.01994662578365207
$ws.Range("D8").Value = 324.334940613

$ws.Range("B9").Value = -135.7697110128234
$ws.Range("C9").Value = 0.08848844495508709
$ws.Range("D9").Value = 238.903778818

$ws.Range("B10").Value = -136.4115378108678
$ws.Range("C10").Value = 0.005124068076634766
$ws.Range("D10").Value = 410.598612087

$ws.Range("B11").Value = -132.4617911649394
$ws.Range("C11").Value = 0.02730893040950624
$ws.Range("D11").Value = 359.950279801
